# [Experimental] Players are no longer fixed to their listed positions.
# Adds an "Is_Pos" column (R) flagging whether a player is shown in a
# position they natively play, and refreshes the player roster data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Is_Pos" header in column R, matching the header style
# already used by the rest of row 1 (bold, bordered, centered).
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(1,18).Value = "Is_Pos"

# Row 2: Jamal Musiala
$ws.Cells.Item(2,1).Value = 486  # ID
$ws.Cells.Item(2,2).Value = "Jamal Musiala"  # Name
$ws.Cells.Item(2,3).Value = "CM"  # Position
$ws.Cells.Item(2,4).Value = 4  # SM
$ws.Cells.Item(2,5).Value = 4  # WF
$ws.Cells.Item(2,6).Value = 81  # Rating
$ws.Cells.Item(2,7).Value = "gold"  # Color
$ws.Cells.Item(2,8).Value = "Bundesliga"  # League
$ws.Cells.Item(2,9).Value = "FC Bayern"  # Club
$ws.Cells.Item(2,10).Value = "Germany"  # Country
$ws.Cells.Item(2,11).Value = 19  # Age
$ws.Cells.Item(2,12).Value = "Rare"  # Rarity
$ws.Cells.Item(2,13).Value = "Right"  # Foot
$ws.Cells.Item(2,14).Value = "Med"  # Attack WR
$ws.Cells.Item(2,15).Value = "Med"  # Defense WR
$ws.Cells.Item(2,16).Value = 700  # Cost
$ws.Cells.Item(2,17).Value = 3  # Chemistry
$ws.Cells.Item(2,18).Value = 1  # Is_Pos

# Row 3: Karim Adeyemi
$ws.Cells.Item(3,1).Value = 514  # ID
$ws.Cells.Item(3,2).Value = "Karim Adeyemi"  # Name
$ws.Cells.Item(3,3).Value = "ST"  # Position
$ws.Cells.Item(3,4).Value = 4  # SM
$ws.Cells.Item(3,5).Value = 3  # WF
$ws.Cells.Item(3,6).Value = 75  # Rating
$ws.Cells.Item(3,7).Value = "gold"  # Color
$ws.Cells.Item(3,8).Value = "Bundesliga"  # League
$ws.Cells.Item(3,9).Value = "Dortmund"  # Club
$ws.Cells.Item(3,10).Value = "Germany"  # Country
$ws.Cells.Item(3,11).Value = 21  # Age
$ws.Cells.Item(3,12).Value = "Rare"  # Rarity
$ws.Cells.Item(3,13).Value = "Left"  # Foot
$ws.Cells.Item(3,14).Value = "High"  # Attack WR
$ws.Cells.Item(3,15).Value = "Med"  # Defense WR
$ws.Cells.Item(3,16).Value = 650  # Cost
$ws.Cells.Item(3,17).Value = 3  # Chemistry
$ws.Cells.Item(3,18).Value = 1  # Is_Pos

# Row 4: Ko Itakura
$ws.Cells.Item(4,1).Value = 542  # ID
$ws.Cells.Item(4,2).Value = "Ko Itakura"  # Name
$ws.Cells.Item(4,3).Value = "CB"  # Position
$ws.Cells.Item(4,4).Value = 3  # SM
$ws.Cells.Item(4,5).Value = 4  # WF
$ws.Cells.Item(4,6).Value = 75  # Rating
$ws.Cells.Item(4,7).Value = "gold"  # Color
$ws.Cells.Item(4,8).Value = "Bundesliga"  # League
$ws.Cells.Item(4,9).Value = "M'gladbach"  # Club
$ws.Cells.Item(4,10).Value = "Japan"  # Country
$ws.Cells.Item(4,11).Value = 25  # Age
$ws.Cells.Item(4,12).Value = "Rare"  # Rarity
$ws.Cells.Item(4,13).Value = "Right"  # Foot
$ws.Cells.Item(4,14).Value = "Med"  # Attack WR
$ws.Cells.Item(4,15).Value = "Med"  # Defense WR
$ws.Cells.Item(4,16).Value = 650  # Cost
$ws.Cells.Item(4,17).Value = 3  # Chemistry
$ws.Cells.Item(4,18).Value = 1  # Is_Pos

# Row 5: Mark Flekken
$ws.Cells.Item(5,1).Value = 562  # ID
$ws.Cells.Item(5,2).Value = "Mark Flekken"  # Name
$ws.Cells.Item(5,3).Value = "GK"  # Position
$ws.Cells.Item(5,4).Value = 1  # SM
$ws.Cells.Item(5,5).Value = 4  # WF
$ws.Cells.Item(5,6).Value = 80  # Rating
$ws.Cells.Item(5,7).Value = "gold"  # Color
$ws.Cells.Item(5,8).Value = "Bundesliga"  # League
$ws.Cells.Item(5,9).Value = "SC Freiburg"  # Club
$ws.Cells.Item(5,10).Value = "Netherlands"  # Country
$ws.Cells.Item(5,11).Value = 29  # Age
$ws.Cells.Item(5,12).Value = "Rare"  # Rarity
$ws.Cells.Item(5,13).Value = "Right"  # Foot
$ws.Cells.Item(5,14).Value = "Med"  # Attack WR
$ws.Cells.Item(5,15).Value = "Med"  # Defense WR
$ws.Cells.Item(5,16).Value = 700  # Cost
$ws.Cells.Item(5,17).Value = 3  # Chemistry
$ws.Cells.Item(5,18).Value = 1  # Is_Pos

# Row 6: Niklas Stark
$ws.Cells.Item(6,1).Value = 812  # ID
$ws.Cells.Item(6,2).Value = "Niklas Stark"  # Name
$ws.Cells.Item(6,3).Value = "CB"  # Position
$ws.Cells.Item(6,4).Value = 2  # SM
$ws.Cells.Item(6,5).Value = 3  # WF
$ws.Cells.Item(6,6).Value = 75  # Rating
$ws.Cells.Item(6,7).Value = "gold"  # Color
$ws.Cells.Item(6,8).Value = "Bundesliga"  # League
$ws.Cells.Item(6,9).Value = "Werder Bremen"  # Club
$ws.Cells.Item(6,10).Value = "Germany"  # Country
$ws.Cells.Item(6,11).Value = 27  # Age
$ws.Cells.Item(6,12).Value = "Rare"  # Rarity
$ws.Cells.Item(6,13).Value = "Right"  # Foot
$ws.Cells.Item(6,14).Value = "Med"  # Attack WR
$ws.Cells.Item(6,15).Value = "Med"  # Defense WR
$ws.Cells.Item(6,16).Value = 650  # Cost
$ws.Cells.Item(6,17).Value = 3  # Chemistry
$ws.Cells.Item(6,18).Value = 1  # Is_Pos

# Row 7: Marc-Oliver Kempf
$ws.Cells.Item(7,1).Value = 2016  # ID
$ws.Cells.Item(7,2).Value = "Marc-Oliver Kempf"  # Name
$ws.Cells.Item(7,3).Value = "CB"  # Position
$ws.Cells.Item(7,4).Value = 2  # SM
$ws.Cells.Item(7,5).Value = 4  # WF
$ws.Cells.Item(7,6).Value = 75  # Rating
$ws.Cells.Item(7,7).Value = "gold"  # Color
$ws.Cells.Item(7,8).Value = "Bundesliga"  # League
$ws.Cells.Item(7,9).Value = "Hertha Berlin"  # Club
$ws.Cells.Item(7,10).Value = "Germany"  # Country
$ws.Cells.Item(7,11).Value = 27  # Age
$ws.Cells.Item(7,12).Value = "Rare"  # Rarity
$ws.Cells.Item(7,13).Value = "Left"  # Foot
$ws.Cells.Item(7,14).Value = "High"  # Attack WR
$ws.Cells.Item(7,15).Value = "Med"  # Defense WR
$ws.Cells.Item(7,16).Value = 650  # Cost
$ws.Cells.Item(7,17).Value = 3  # Chemistry
$ws.Cells.Item(7,18).Value = 1  # Is_Pos

# Row 8: Ridle Baku
$ws.Cells.Item(8,1).Value = 2101  # ID
$ws.Cells.Item(8,2).Value = "Ridle Baku"  # Name
$ws.Cells.Item(8,3).Value = "RM"  # Position
$ws.Cells.Item(8,4).Value = 3  # SM
$ws.Cells.Item(8,5).Value = 4  # WF
$ws.Cells.Item(8,6).Value = 80  # Rating
$ws.Cells.Item(8,7).Value = "gold"  # Color
$ws.Cells.Item(8,8).Value = "Bundesliga"  # League
$ws.Cells.Item(8,9).Value = "VfL Wolfsburg"  # Club
$ws.Cells.Item(8,10).Value = "Germany"  # Country
$ws.Cells.Item(8,11).Value = 24  # Age
$ws.Cells.Item(8,12).Value = "Rare"  # Rarity
$ws.Cells.Item(8,13).Value = "Right"  # Foot
$ws.Cells.Item(8,14).Value = "High"  # Attack WR
$ws.Cells.Item(8,15).Value = "Med"  # Defense WR
$ws.Cells.Item(8,16).Value = 700  # Cost
$ws.Cells.Item(8,17).Value = 3  # Chemistry
$ws.Cells.Item(8,18).Value = 1  # Is_Pos

# Row 9: Josip Šutalo
$ws.Cells.Item(9,1).Value = 2372  # ID
$ws.Cells.Item(9,2).Value = "Josip Šutalo"  # Name
$ws.Cells.Item(9,3).Value = "CB"  # Position
$ws.Cells.Item(9,4).Value = 2  # SM
$ws.Cells.Item(9,5).Value = 3  # WF
$ws.Cells.Item(9,6).Value = 75  # Rating
$ws.Cells.Item(9,7).Value = "gold"  # Color
$ws.Cells.Item(9,8).Value = "Liga Hrvatska (CRO 1)"  # League
$ws.Cells.Item(9,9).Value = "Dinamo Zagreb"  # Club
$ws.Cells.Item(9,10).Value = "Croatia"  # Country
$ws.Cells.Item(9,11).Value = 22  # Age
$ws.Cells.Item(9,12).Value = "Rare"  # Rarity
$ws.Cells.Item(9,13).Value = "Right"  # Foot
$ws.Cells.Item(9,14).Value = "Med"  # Attack WR
$ws.Cells.Item(9,15).Value = "High"  # Defense WR
$ws.Cells.Item(9,16).Value = 650  # Cost
$ws.Cells.Item(9,17).Value = 0  # Chemistry
$ws.Cells.Item(9,18).Value = 0  # Is_Pos

# Row 10: Rico Henry
$ws.Cells.Item(10,1).Value = 8198  # ID
$ws.Cells.Item(10,2).Value = "Rico Henry"  # Name
$ws.Cells.Item(10,3).Value = "LWB"  # Position
$ws.Cells.Item(10,4).Value = 3  # SM
$ws.Cells.Item(10,5).Value = 2  # WF
$ws.Cells.Item(10,6).Value = 75  # Rating
$ws.Cells.Item(10,7).Value = "gold"  # Color
$ws.Cells.Item(10,8).Value = "Premier League"  # League
$ws.Cells.Item(10,9).Value = "Brentford"  # Club
$ws.Cells.Item(10,10).Value = "England"  # Country
$ws.Cells.Item(10,11).Value = 25  # Age
$ws.Cells.Item(10,12).Value = "Rare"  # Rarity
$ws.Cells.Item(10,13).Value = "Left"  # Foot
$ws.Cells.Item(10,14).Value = "High"  # Attack WR
$ws.Cells.Item(10,15).Value = "High"  # Defense WR
$ws.Cells.Item(10,16).Value = 650  # Cost
$ws.Cells.Item(10,17).Value = 0  # Chemistry
$ws.Cells.Item(10,18).Value = 0  # Is_Pos

# Row 11: Alfonso Espino
$ws.Cells.Item(11,1).Value = 8821  # ID
$ws.Cells.Item(11,2).Value = "Alfonso Espino"  # Name
$ws.Cells.Item(11,3).Value = "LB"  # Position
$ws.Cells.Item(11,4).Value = 3  # SM
$ws.Cells.Item(11,5).Value = 3  # WF
$ws.Cells.Item(11,6).Value = 78  # Rating
$ws.Cells.Item(11,7).Value = "gold"  # Color
$ws.Cells.Item(11,8).Value = "LaLiga Santander"  # League
$ws.Cells.Item(11,9).Value = "Cádiz CF"  # Club
$ws.Cells.Item(11,10).Value = "Uruguay"  # Country
$ws.Cells.Item(11,11).Value = 31  # Age
$ws.Cells.Item(11,12).Value = "Rare"  # Rarity
$ws.Cells.Item(11,13).Value = "Left"  # Foot
$ws.Cells.Item(11,14).Value = "High"  # Attack WR
$ws.Cells.Item(11,15).Value = "Med"  # Defense WR
$ws.Cells.Item(11,16).Value = 600  # Cost
$ws.Cells.Item(11,17).Value = 0  # Chemistry
$ws.Cells.Item(11,18).Value = 0  # Is_Pos

# Row 12: Ruben Vargas
$ws.Cells.Item(12,1).Value = 9424  # ID
$ws.Cells.Item(12,2).Value = "Ruben Vargas"  # Name
$ws.Cells.Item(12,3).Value = "LM"  # Position
$ws.Cells.Item(12,4).Value = 4  # SM
$ws.Cells.Item(12,5).Value = 4  # WF
$ws.Cells.Item(12,6).Value = 75  # Rating
$ws.Cells.Item(12,7).Value = "gold"  # Color
$ws.Cells.Item(12,8).Value = "Bundesliga"  # League
$ws.Cells.Item(12,9).Value = " FC Augsburg"  # Club
$ws.Cells.Item(12,10).Value = "Switzerland"  # Country
$ws.Cells.Item(12,11).Value = 24  # Age
$ws.Cells.Item(12,12).Value = "Rare"  # Rarity
$ws.Cells.Item(12,13).Value = "Right"  # Foot
$ws.Cells.Item(12,14).Value = "High"  # Attack WR
$ws.Cells.Item(12,15).Value = "Med"  # Defense WR
$ws.Cells.Item(12,16).Value = 650  # Cost
$ws.Cells.Item(12,17).Value = 3  # Chemistry
$ws.Cells.Item(12,18).Value = 1  # Is_Pos
